$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.130.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2822"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.859.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.117"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6686"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "280.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.153.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.479"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.100.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007244"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.133"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.314"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.902"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.344"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09599"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.397"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.468"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.096"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04649"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7003"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.092"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.528"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8517"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.918"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4159"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "988.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.164"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1138"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.83%  "
